$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 6, shifting existing rows 6-10 down to 7-11
$ws.Rows.Item(6).Insert()

# Fill in the new row 6 with data
$ws.Range("A6").Value = 7
$ws.Range("B6").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C6").Value = "Ñuble"
$ws.Range("D6").NumberFormat = $ws.Range("D7").NumberFormat
$ws.Range("D6").Value = 45036
$ws.Range("E6").Value = 16
$ws.Range("F6").Value = "Fruta"
$ws.Range("G6").Value = 100107
$ws.Range("H6").Value = "Otros"
$ws.Range("I6").Value = 100107011
$ws.Range("J6").Value = "Tuna"
$ws.Range("K6").Value = "Sin especificar"
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 60
$ws.Range("N6").Value = 15000
$ws.Range("O6").Value = 16000
$ws.Range("P6").Value = 15500
$ws.Range("Q6").Value = "$/caja 18 kilos"
$ws.Range("R6").Value = "Región Metropolitana"
$ws.Range("S6").Value = 861
$ws.Range("T6").Value = 18
